# fix: cleaning script up after debugging
#
# Updates the "Recipes" sheet:
#  - rows 2 & 3 change Meal Type from "Lunch" to "Dinner" and clear the
#    stray "N/A" Calories placeholder
#  - five new recipe rows (4-8) are appended

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Recipes")

function Set-EmptyCell($range) {
    # Assigning an empty string clears the cell, but we still want the
    # cell itself to remain present in the sheet (just with no value).
    # Touching a formatting property that is already at its default
    # value keeps the cell "alive" without altering its appearance.
    $range.Value = ""
    $range.Font.Bold = $False
}

# --- Update existing rows 2 and 3 ---
$ws.Range("B2").Value = "Dinner"
Set-EmptyCell $ws.Range("F2")

$ws.Range("B3").Value = "Dinner"
Set-EmptyCell $ws.Range("F3")

# --- Add new row 4 ---
$ws.Range("A4").Value = "Asparagus and Pea Soup: Real Convenience Food"
$ws.Range("B4").Value = "Breakfastr"
Set-EmptyCell $ws.Range("C4")
$ws.Range("D4").Value = "https://spoonacular.com/recipes/716406"
$ws.Range("E4").Value = "1 bag of frozen organic asparagus (preferably thawed), 1T EVOO (extra virgin olive oil), a couple of garlic cloves, 1/2 onion, 2-3c of frozen organic peas, 1 box low-sodium vegetable broth"
Set-EmptyCell $ws.Range("F4")

# --- Add new row 5 ---
$ws.Range("A5").Value = "Garlicky Kale"
$ws.Range("B5").Value = "Breakfastr"
Set-EmptyCell $ws.Range("C5")
$ws.Range("D5").Value = "https://spoonacular.com/recipes/644387"
$ws.Range("E5").Value = "3 tablespoons balsamic vinegar, 1 clove garlic, minced, 1 bunch curly kale, stems removed and chopped, Olive oil"
Set-EmptyCell $ws.Range("F5")

# --- Add new row 6 ---
$ws.Range("A6").Value = "Powerhouse Almond Matcha Superfood Smoothie"
$ws.Range("B6").Value = "Breakfast"
Set-EmptyCell $ws.Range("C6")
$ws.Range("D6").Value = "https://spoonacular.com/recipes/756814"
$ws.Range("E6").Value = "2 tablespoons unsalted natural almond butter, 1 1/2 cups unsweetened almond milk, 1 medium frozen banana, 2 teaspoons chia seeds, 1 cup baby kale, packed, 1/2 cup frozen mango pieces, 1 tablespoon matcha green tea powder, 3/4 cup frozen pineapple, 1/2 teaspoon vanilla extract"
Set-EmptyCell $ws.Range("F6")

# --- Add new row 7 ---
$ws.Range("A7").Value = "Butternut Squash Frittata"
$ws.Range("B7").Value = "Breakfast"
Set-EmptyCell $ws.Range("C7")
$ws.Range("D7").Value = "https://spoonacular.com/recipes/636589"
$ws.Range("E7").Value = "1 large butternut squash, peeled, seeded, thinly sliced (with a mandoline), 1/2 oz goat cheese, 1/2 cup liquid egg substitute, 2 tbsp. non-fat milk, Pepper to taste"
Set-EmptyCell $ws.Range("F7")

# --- Add new row 8 ---
$ws.Range("A8").Value = "Doughnuts"
$ws.Range("B8").Value = "Breakfast"
Set-EmptyCell $ws.Range("C8")
$ws.Range("D8").Value = "https://spoonacular.com/recipes/716276"
$ws.Range("E8").Value = "1.5 cups of flour, 30 ml honey, 1 tablespoon of powdered milk, 1/2 teaspoon salt, 150 ml warm water, 1 teaspoon yeast"
Set-EmptyCell $ws.Range("F8")
